# feat: Find instances with all disagreeing annotators
#
# The template's example row (row 4) previously had a sample/default
# value of 0 in the "number of disagreeing annotators" column (C4),
# left over from when the template was authored. Clear it out so the
# template ships blank, and move the active selection to C4 (the cell
# that actually drives the "find instances with all disagreeing
# annotators" scenario) instead of the stray A9 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover sample value in C4 (was <v>0</v>), keeping its
# style/data validation intact.
$ws.Range("C4").ClearContents()

# Update the active selection/cell to C4.
$ws.Range("C4").Select()
